$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.144.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.282.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.40%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "155.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15,425.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "94.61"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.79%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.492"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "35.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +7.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0805"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "

# Row 13
$ws.Range("E13").Value = "  -2.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.636.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.26%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.280.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.795"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.20%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.064.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0919"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "244.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.05%  "

# Row 26
$ws.Range("E26").Value = "  +0.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "

# Row 31
$ws.Range("E31").Value = "  +1.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.11%  "

# Row 34
$ws.Range("E34").Value = "  -0.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0754"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "

# Row 36
$ws.Range("E36").Value = "  +1.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.55%  "

# Row 39
$ws.Range("E39").Value = "  -0.25%  "

# Row 40
$ws.Range("E40").Value = "  -0.18%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.021.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.39%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.42%  "

# Row 46
$ws.Range("E46").Value = "  +1.75%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.37%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.73%  "
